# Insert a new weekly record row before existing row 102, shifting rows
# 102:195 down to 103:196, then populate the new row 102 with a duplicate
# of the original row 102 data but with an updated date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 102:195 down by one to make room for the new record. Excel
# copies the formatting of the row above into the newly created row.
$ws.Rows.Item(102).Insert()

# Row 103 now holds what used to be row 102's data (prior to the shift).
# Duplicate that data into the newly blank row 102.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(102, $col).Value = $ws.Cells.Item(103, $col).Value2
}

# The new record uses an updated date (2022-04-18 => serial 44669) while
# every other field is carried over unchanged from the duplicated row.
$ws.Cells.Item(102, 4).Value = 44669

$wb.Save()
